$wb = $excel.ActiveWorkbook

# Duplicate the "Requests" sheet (it already carries the bold/bordered
# header style we need) and move the copy to the end of the workbook so
# we don't have to rebuild the style table by hand.
$src = $wb.Worksheets.Item("Requests")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)

# The copy is now the last sheet - turn it into "Server Use" with its own
# header row: Event | BS | Use
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Server Use"
$new.Range("B1").Value = "Event"
$new.Range("C1").Value = "BS"
$new.Range("D1").Value = "Use"
$new.Range("E1:G1").Clear()

# Restore "Requests" as the active sheet, matching the original workbook.
$wb.Worksheets.Item("Requests").Activate()
